# Fix Training Data Issue (#48)
# The "Date" column (BF) was off by a day because of how NBA stats were
# originally scraped/shown. Correct every date value in column BF
# (rows 2-31) from the old "6-10-2011-12" label to the proper ISO date
# "2012-06-10".
#
# Note: assigning a date-shaped string straight to .Value/.Value2/.Formula
# (or via Range.Replace) causes Excel to auto-recognize it as a date and
# store it as a numeric date serial (with a new number-format style).
# To keep the corrected value as plain text (matching the rest of the
# sheet, which stores this column as text), we write it as a literal
# text formula first and then convert the range to a static value via
# Copy / Paste-Special-Values, which preserves the "text" cell type
# without adding any new number formats/styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCol = $ws.Range("BF2:BF31")

$dateCol.Formula = "=""2012-06-10"""
$dateCol.Copy() | Out-Null
$dateCol.PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = $false
